$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: remove all existing values from the old invoice table
$ws.Cells.ClearContents()

# ---- Row 1: header labels ----
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "DOB"
$ws.Range("C1").Value = "Nationality"
$ws.Range("D1").Value = "Place of Issue"
$ws.Range("E1").Value = "Date of Issue"
$ws.Range("F1").Value = "Date of Expiry"

# Give the two new header cells (E1, F1) the same look as the rest of the
# header row (bold font, border, centered) by copying the format from D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null

# ---- Data rows, filled column by column (matches extraction order of the
# multi-model document processing pipeline: each field/column is resolved
# for every row before moving to the next field) ----

# Column A: Name / Date of Birth (label)
$ws.Range("A2").Value = "Suresh Sharma"
$ws.Range("A3").Value = "Date of Birth"

# Column B: DOB value
# "01/10/1999" looks like a date to Excel's auto-detection. Enter it as a
# literal-text formula instead (so it is never date-parsed), then collapse
# the formula to a plain value in place with Copy/PasteSpecial values-only.
# This keeps the cell's style untouched (no extra number-format style).
$ws.Range("B2").Formula = '="01/10/1999"'
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4163) | Out-Null

# Column C: Nationality / Passport (label)
$ws.Range("C2").Value = "Indian"
$ws.Range("C3").Value = "Passport"

# Column D: Place of Issue / Date of Issue (label, reuses header text)
$ws.Range("D2").Value = "Hyderabad"
$ws.Range("D3").Value = "Date of Issue"

# Column E: Date of Issue value
$ws.Range("E2").Value = "14/08/2023"

# Column F: Date of Expiry value (also looks like a date, same guard)
$ws.Range("F2").Formula = '="11/08/2033"'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4163) | Out-Null
